$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # Some "Price" values in column D are plain decimal numbers (e.g. "211.77").
    # A normal .Value assignment lets Excel's type inference turn those into
    # numeric cells, but the source sheet keeps every Price cell as text.
    # Forcing the number format to Text before the write keeps the cell a
    # string; ClearFormats() afterwards drops the now-unneeded explicit
    # number format so the cell's style index goes back to the sheet default.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "27.915.16"
$ws.Range("E2").Value = "  -0.24%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.632.65"
$ws.Range("E3").Value = "  -0.88%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
Set-TextValue "D5" "211.77"
$ws.Range("E5").Value = "  -0.79%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.94%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.08%  "

# Row 8 - Solana
Set-TextValue "D8" "23.27"
$ws.Range("E8").Value = "  -0.77%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -2.91%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.20%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0881"
$ws.Range("E11").Value = "  +1.04%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.864.70"
$ws.Range("E12").Value = "  -0.88%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.634.07"
$ws.Range("E13").Value = "  -0.99%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.17%  "

# Row 16 - Litecoin
Set-TextValue "D16" "65.25"
$ws.Range("E16").Value = "  -0.51%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "27.917.31"
$ws.Range("E17").Value = "  -0.16%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "230.53"
$ws.Range("E18").Value = "  -1.05%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  -0.14%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -2.59%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.03%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.71%  "

# Row 23 - Avalanche
Set-TextValue "D23" "10.36"
$ws.Range("E23").Value = "  -2.88%  "

# Row 24 - Toncoin
Set-TextValue "D24" "2.06"
$ws.Range("E24").Value = "  -3.72%  "

# Row 25 - Monero
Set-TextValue "D25" "154.39"
$ws.Range("E25").Value = "  +1.21%  "

# Row 26 - Cosmos
Set-TextValue "D26" "6.97"
$ws.Range("E26").Value = "  +0.63%  "

# Row 27 - now EthereumClassic (was Stellar)
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D27" "15.65"
$ws.Range("E27").Value = "  -0.73%  "

# Row 28 - now Stellar (was EthereumClassic)
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D28" "0.111"
$ws.Range("E28").Value = "  -0.67%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  -0.15%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -1.06%  "

# Row 31 - Hedera
Set-TextValue "D31" "0.0483"
$ws.Range("E31").Value = "  -0.02%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.89%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -0.58%  "

# Row 34 - Maker
Set-TextValue "D34" "1.401.31"
$ws.Range("E34").Value = "  -3.02%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +0.13%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  +9.68%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  +1.41%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +0.40%  "

# Row 39 - ImmutableX
Set-TextValue "D39" "0.562"
$ws.Range("E39").Value = "  +0.24%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -1.96%  "

# Row 41 - WEMIXToken
$ws.Range("E41").Value = "  -0.23%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  -0.11%  "

# Row 43 - Aave
Set-TextValue "D43" "66.80"
$ws.Range("E43").Value = "  -3.81%  "

# Row 44 - FraxShare
$ws.Range("E44").Value = "  +2.67%  "

# Row 45 - RenderToken
$ws.Range("E45").Value = "  +1.10%  "

# Row 46 - MXToken
$ws.Range("E46").Value = "  -1.25%  "

# Row 47 - RocketPoolETH
Set-TextValue "D47" "1.774.27"
$ws.Range("E47").Value = "  -0.93%  "

# Row 48 - Quant
Set-TextValue "D48" "87.66"
$ws.Range("E48").Value = "  -1.45%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  +0.43%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  -0.77%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  -0.35%  "
